# daily auto push: 2025-10-13 07:29 UTC
# Append the new day's row (row 99) to the bottom of the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as plain text so the date-like string
# "2025/10/13" is stored verbatim (matching the rest of the column)
# instead of being auto-converted into a date serial number.
$ws.Range("A99").NumberFormat = "@"
$ws.Range("A99").Value = "2025/10/13"
$ws.Range("B99").Value = "月"
$ws.Range("C99").Value = 16
$ws.Range("D99").Value = 19

# Drop the temporary text-format override so the new row's cells end up
# with the same (default) style as the rest of the data rows.
$ws.Range("A99").ClearFormats()
